$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.636641502380371
$ws.Range("B1").Value = 3.650352954864502
$ws.Range("C1").Value = 2.107701778411865
$ws.Range("D1").Value = 1.716450810432434
$ws.Range("E1").Value = 1.608127474784851
